$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.282.65"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.13%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.868.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.42%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7241"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.38%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.41%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07834"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.73%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3091"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.69%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.26"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.43%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08251"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.26%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.875.43"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.08%  "

$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7227"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.30%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.244"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.74%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "90.85"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.87%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.266.07"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.02%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.859"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.44%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "244.09"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.92%  "

$ws.Range("E19").Value = "  +0.72%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.52%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.102.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.81%  "

$ws.Range("E22").Value = "  -0.05%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.969"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.19%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.04%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1598"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +11.94%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.91"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.05%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.972"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.07%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.25"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.26%  "

$ws.Range("E29").Value = "  -1.67%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.496"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.36%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.401"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.74%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.103"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.99%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05198"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.00%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.934"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.12%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.187"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.72%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7291"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.66%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01858"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.86%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.698"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.32%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.175.54"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.07%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9041"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.04%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.103"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.84%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "72.55"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.47%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.02%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.86"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.60%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5281"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.36%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.999.77"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.79%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.783"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.20%  "

$ws.Range("B49").Value = "SynthetixNetwork"
$ws.Range("C49").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.893"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.87%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.299"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.97%  "

$ws.Range("B51").Value = "TheSandbox"
$ws.Range("C51").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4271"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.96%  "
